# Edit script: turns the "Unraveling the Matrix of Time" essay into the
# "Exploring the Vast World of Mathematics" essay, per the commit diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: not found -> $old"
    }
}

# --- Title / byline / contact ---------------------------------------------
Replace-Text "Unraveling the Matrix of Time: A Quest for Order" "Exploring the Vast World of Mathematics: A Journey Through Numbers and Patterns"
Replace-Text "Jasmine Davis" "Dr Deborah Fowler"
Replace-Text "jasmine" "fowler"
Replace-Text "davis@astralnexus" "deborah@hillcrestacademy"

# --- Body paragraph 1 (Intro) ----------------------------------------------
Replace-Text "As humans, we are perpetually bound to the construct of time; its relentless progression shapes our perception of reality and influences the very fabric of our existence" "Mathematics, the science of numbers and patterns, has captivated the minds of thinkers throughout history"

Replace-Text " Yet, concealed within the enigma of temporal flow lies a symphony of patterns, an underlying order that hums beneath the clamor of chaos" " It is a language of its own, describing the underlying principles that govern the universe"

Replace-Text " This exploration delves into the nature of time, investigating its paradoxical essence and unveiling the secret harmonies woven into its tapestry" " From the intricate designs of nature to the technological marvels of the modern world, mathematics plays a pivotal role in shaping our understanding of the cosmos"

# Two runs removed (" We will traverse..." + ".") and the following run's text replaced.
Replace-Text " We will traverse diverse disciplines, from the cosmic elegance of physics to the intricate machinery of computer science, in pursuit of unveiling the matrix of time, its elusive enigma. The journey promises revelations that transcend our conventional understanding, inviting us to perceive time as both an absolute entity and a fluid construct" " As we delve into the realm of mathematics, we embark on a journey of discovery, unraveling the secrets of numbers and patterns, and gaining a deeper appreciation for the beauty and power of this remarkable subject"

# --- Body paragraph 1 (Physics section -> Math usefulness section) --------
Replace-Text "Venturing into the realm of physics, General Relativity, the monumental work of Albert Einstein, unveils a universe where time and space intertwine, forming an indissoluble fabric" "Mathematics, present in every aspect of life, offers a gateway to understanding the world around us"

Replace-Text " Gravity, the invisible force that orchestrates the cosmic ballet, bends and warps this spacetime continuum, inducing temporal distortions and gravitational time dilation" " It allows us to quantify, measure, and analyze data, making sense of complex information"

Replace-Text " We encounter a realm where the ticking of time varies, a symphony of rhythms under the baton of celestial bodies" " From navigating through the urban grid to predicting weather patterns, mathematics provides the tools essential for solving problems and making informed decisions"

# Two runs removed (" The nature of time..." + ".") and the following run's text replaced.
Replace-Text " The nature of time becomes pliable, subject to the influence of mass and energy, challenging our intuitive grasp of its uniform flow. Furthermore, as we traverse the vast spectrum of light, a peculiar phenomenon emerges - the Twins Paradox, a conundrum that pits relativity against common sense" " It is a key component of various fields, including science, engineering, economics, and even art"

# --- Body paragraph 1 (Computer science -> creativity section) ------------
Replace-Text "Descending from the cosmic to the digital realm, we encounter computer science, a discipline that has unlocked the secrets of virtual time" "Furthermore, mathematics stimulates creativity and critical thinking"

Replace-Text " Here, time becomes quantized, reduced to discrete increments, akin to a choreographed dance of transistors and logic gates" " By engaging in mathematical puzzles and problems, we are challenged to think outside the box, explore new approaches, and develop logical reasoning skills"

Replace-Text " Algorithms, intricate sequences of instructions, slice and dice the continuum, creating artificial time frames within which computation unfolds" " Mathematics encourages us to question, analyze, and seek patterns, nurturing an inquisitive mindset that is crucial for lifelong learning and success"

# Three runs removed (" This digital realm..." + "." + " " + lastRenderedPageBreak/"Through these explorations...")
# and the text of the first of those runs replaced; lastRenderedPageBreak marker is dropped.
Replace-Text " This digital realm serves as a laboratory for exploring alternative temporal structures, allowing us to manipulate and distort time in ways that would confound our physical counterparts. Through these explorations, we glean insights into the malleability of time, its susceptibility to human design" " As we delve deeper into the mathematical landscape, we discover that it is not merely a collection of abstract concepts"

# --- Summary paragraph -------------------------------------------------
Replace-Text "Our quest through the tapestry of time reveals its elusive and paradoxical nature, where absoluteness and fluidity coexist" "Exploring the vast world of mathematics is an enriching journey that uncovers the beauty and power of numbers and patterns"

Replace-Text " From the symphony of patterns in General Relativity to the quantized rhythms of computer science, we catch glimpses of time's hidden order" " Mathematics is not just a subject taught in school; it's a language that describes the principles governing the universe"

Replace-Text " Yet, the essence of time, its ultimate enigma, remains shrouded in mystery" " It touches every aspect of life, helping us understand the world around us"

Replace-Text " Perhaps, like the cosmic dancer Shiva, time performs its cosmic dance, its secrets forever concealed within its ceaseless motion" " Mathematics enhances critical thinking, creativity, and problem-solving skills, which are essential for success in various fields and lifelong learning"

# New sentence appended at the end of the Summary paragraph, right before the
# final "." run (which is left untouched).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$find = $d.Content
$find.Find.Execute("Mathematics enhances critical thinking, creativity, and problem-solving skills, which are essential for success in various fields and lifelong learning") | Out-Null
$insertPoint = $d.Range($find.End, $find.End)
$insertPoint.InsertAfter(". Join us as we explore the wonders of mathematics, unraveling its secrets and unlocking its potential to transform our understanding of the cosmos")

# --- New empty paragraph at the very end of the document -------------------
$d.Content.InsertParagraphAfter()
